$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for rows 2-13 (4 senders x 3 targets; "Resolving-Mac" dropped as a target cluster)
$data = @(
    ,@("ECs", "Efnb2", "Ephb6", "ECs", 3, 1, 45.76217133333333, 137.286514, 0.6763939203605134, 0.6763939203605135, 2, 0.6666666666666666, 0.168144, 0.504432, 0.05446245276675245, 0.05446245276675245, 7.694634536672, 69.251710830048, 0.03683807193935298, 0.03683807193935298)
    ,@("ECs", "Efnb2", "Ephb6", "FAPs", 3, 1, 45.76217133333333, 137.286514, 0.6763939203605134, 0.6763939203605135, 3, 1, 1.279928333333333, 3.839785, 0.4145734394268892, 0.4145734394268892, 58.57229968438778, 527.15069715949, 0.2804149539712955, 0.2804149539712955)
    ,@("ECs", "Efnb2", "Ephb6", "MuSCs", 3, 1, 45.76217133333333, 137.286514, 0.6763939203605134, 0.6763939203605135, 3, 1, 1.639265666666667, 4.917797, 0.5309641078063584, 0.5309641078063583, 75.0163562988509, 675.1472066896581, 0.359140894449865, 0.359140894449865)
    ,@("FAPs", "Efnb2", "Ephb6", "ECs", 3, 1, 9.395935333333332, 28.187806, 0.1388778842960613, 0.1388778842960613, 2, 0.6666666666666666, 0.168144, 0.504432, 0.05446245276675245, 0.05446245276675245, 1.579870150688, 14.218831356192, 0.007563630213820751, 0.007563630213820752)
    ,@("FAPs", "Efnb2", "Ephb6", "FAPs", 3, 1, 9.395935333333332, 28.187806, 0.1388778842960613, 0.1388778842960613, 3, 1, 1.279928333333333, 3.839785, 0.4145734394268892, 0.4145734394268892, 12.02612385130111, 108.23511466171, 0.0575750821529477, 0.0575750821529477)
    ,@("FAPs", "Efnb2", "Ephb6", "MuSCs", 3, 1, 9.395935333333332, 28.187806, 0.1388778842960613, 0.1388778842960613, 3, 1, 1.639265666666667, 4.917797, 0.5309641078063584, 0.5309641078063583, 15.40243419815355, 138.621907783382, 0.07373917192929286, 0.07373917192929286)
    ,@("MuSCs", "Efnb2", "Ephb6", "ECs", 3, 1, 12.29750866666667, 36.892526, 0.1817649787009828, 0.1817649787009828, 2, 0.6666666666666666, 0.168144, 0.504432, 0.05446245276675245, 0.05446245276675245, 2.067752297248, 18.609770675232, 0.00989936656715204, 0.00989936656715204)
    ,@("MuSCs", "Efnb2", "Ephb6", "FAPs", 3, 1, 12.29750866666667, 36.892526, 0.1817649787009828, 0.1817649787009828, 3, 1, 1.279928333333333, 3.839785, 0.4145734394268892, 0.4145734394268892, 15.73992977187889, 141.65936794691, 0.07535493238742169, 0.07535493238742168)
    ,@("MuSCs", "Efnb2", "Ephb6", "MuSCs", 3, 1, 12.29750866666667, 36.892526, 0.1817649787009828, 0.1817649787009828, 3, 1, 1.639265666666667, 4.917797, 0.5309641078063584, 0.5309641078063583, 20.15888374280244, 181.429953685222, 0.09651067974640905, 0.09651067974640903)
    ,@("Resolving-Mac", "Efnb2", "Ephb6", "ECs", 3, 1, 0.2004796666666666, 0.601439, 0.002963216642442438, 0.002963216642442439, 2, 0.6666666666666666, 0.168144, 0.504432, 0.05446245276675245, 0.05446245276675245, 0.03370945307199999, 0.303385077648, 0.0001613840464266761, 0.0001613840464266761)
    ,@("Resolving-Mac", "Efnb2", "Ephb6", "FAPs", 3, 1, 0.2004796666666666, 0.601439, 0.002963216642442438, 0.002963216642442439, 3, 1, 1.279928333333333, 3.839785, 0.4145734394268892, 0.4145734394268892, 0.2565996056238888, 2.309396450615, 0.00122847091522436, 0.00122847091522436)
    ,@("Resolving-Mac", "Efnb2", "Ephb6", "MuSCs", 3, 1, 0.2004796666666666, 0.601439, 0.002963216642442438, 0.002963216642442439, 3, 1, 1.639265666666667, 4.917797, 0.5309641078063584, 0.5309641078063583, 0.3286394344314444, 2.957754909883, 0.001573361680791402, 0.001573361680791402)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowVals[$c]
    }
}

# Old sheet had 16 data rows (4x4 grid); new data only needs 12 rows (4x3 grid) -> delete the trailing rows
$ws.Range("A14:T17").EntireRow.Delete()

Write-Host "applied"
